$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "content header" row (row 18): column labels for the line
# items table below the invoice address block.
$ws.Range("B18").Value = "Désignation"
$ws.Range("F18").Value = "Quantité"
$ws.Range("G18").Value = "P.U. TVAC"
$ws.Range("H18").Value = "Total TVAC"
$ws.Range("I18").Value = "Code TVA"

# Give the new header cells the same (plain Calibri 11) formatting already
# used by the other standalone text cells (e.g. "Doit pour vente"), by
# copying that cell's format instead of touching Font.* directly, so we
# reuse the existing style index rather than minting a new font/style.
$ws.Range("E13").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("E13").Copy()
$ws.Range("F18").PasteSpecial(-4122)
$ws.Range("E13").Copy()
$ws.Range("G18").PasteSpecial(-4122)
$ws.Range("E13").Copy()
$ws.Range("H18").PasteSpecial(-4122)
$ws.Range("E13").Copy()
$ws.Range("I18").PasteSpecial(-4122)

# Row 18 is taller than the default to make room for wrapped header text.
$ws.Rows.Item(18).RowHeight = 30

# Column H (8) gets an explicit width so the "Total TVAC" header fits.
$ws.Columns.Item(8).ColumnWidth = 10.166666666666666
